$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.795.40"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "2.527.85"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.57"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.82"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.577"
$ws.Range("E7").Value = "  -1.61%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.18"
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0809"
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("E13").Value = "  -3.28%  "
$ws.Range("D14").Value = "2.915.92"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.30"
$ws.Range("E15").Value = "  -2.86%  "
$ws.Range("D16").Value = "2.505.36"
$ws.Range("E16").Value = "  -0.74%  "
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "42.866.30"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.76"
$ws.Range("E19").Value = "  +4.70%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.87"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("E22").Value = "  -2.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.03"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.95"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.67"
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("E28").Value = "  +4.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.85"
$ws.Range("E29").Value = "  +8.32%  "
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.41"
$ws.Range("E32").Value = "  +2.46%  "
$ws.Range("E33").Value = "  +4.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.43"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("E35").Value = "  +2.77%  "
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0779"
$ws.Range("E37").Value = "  -0.99%  "
$ws.Range("E38").Value = "  -1.80%  "
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.51"
$ws.Range("E40").Value = "  -4.43%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.34"
$ws.Range("E41").Value = "  +15.21%  "
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.20%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.33"
$ws.Range("E45").Value = "  -1.76%  "
$ws.Range("D46").Value = "2.053.90"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.33"
$ws.Range("E47").Value = "  +1.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "107.80"
$ws.Range("E48").Value = "  +6.24%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.93"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.56"
$ws.Range("E50").Value = "  +3.90%  "
$ws.Range("D51").Value = "2.767.61"
$ws.Range("E51").Value = "  +0.41%  "
